$d = $word.ActiveDocument

# The document currently ends with an empty, numbered (ListParagraph) bullet.
# Fill that paragraph with the first new bullet's text, then add a second
# new bullet paragraph (same list/style) after it with the second text.

$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastRange.InsertAfter("Decide two best depths for each system, and have Kun carry out a few more simulations")

$newPara = $d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertAfter("Look at earlier simulations where Tat position was not constrained and compare against exp.")
